$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per the diff
$ws.Range("B2").Value = 3
$ws.Range("D5").Value = 15
$ws.Range("D6").Value = 45
$ws.Range("D7").Value = 32
$ws.Range("D9").Value = 41

# Update the active selection on the sheet (was G16 -> now D5)
$ws.Range("D5").Select()
